# "New Test case added" - adds a 4th Framework_002/iMacs test-case row
# (Suvam / 25, LimeSquare, City Road) cloned from row 3, switches the
# Browser column from Mozilla to Chrome on every data row, and drops the
# stale Result (Pass/Pass) values from column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clone row 3 (style + values) down into row 4 for the new test case.
$ws.Range("A3:M3").Copy($ws.Range("A4:M4"))

# 2. The old "Result" column values are no longer populated; only the
#    header (N1) remains.
$ws.Range("N2:N3").ClearContents()

# 3. Differentiate the new row's First Name / Address.
$ws.Range("G4").Value2 = "Suvam"
$ws.Range("I4").Value2 = "25, LimeSquare, City Road"

# 4. Switch Browser from Mozilla to Chrome for every test case.
$ws.Range("D2").Value2 = "Chrome"
$ws.Range("D3").Value2 = "Chrome"
$ws.Range("D4").Value2 = "Chrome"

# 5. Wire up the Password/Email hyperlinks for the new row.
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("M4"), "mailto:tooolsqa@gmail.com")

# Re-apply the bordered cell formatting that Hyperlinks.Add overwrote with
# the generic "Hyperlink" style.
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("M3").Copy()
$ws.Range("M4").PasteSpecial(-4122)

# 6. Extend the Product Type / Product Number validation lists to row 4.
$ws.Range("E2:E4").Validation.Delete()
$ws.Range("E2:E4").Validation.Add(3, 1, 1, '"Accessories, iMacs, iPads, iPhones"')
$ws.Range("F2:F4").Validation.Delete()
$ws.Range("F2:F4").Validation.Add(3, 1, 1, '"Product 1, Product 2, Product 3, Product 4"')

# 7. Match the author's final selection.
$ws.Range("F10").Select() | Out-Null
